$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "last report" marker cell (written first so its shared string
#     lands at the same index the source workbook ends up with) ---
$ws.Range("I32").Value = " "

# --- Header date label ---
$ws.Range("B1").Value = "29.01.2025"

# --- Credit section date/payment label ---
$ws.Range("F34").Value = "30.01.2025 payment "

# --- Stock quantity updates (row 9-17) ---
$ws.Range("C9").Value = 342876
$ws.Range("C10").Value = 190
$ws.Range("C11").Value = 590
$ws.Range("C12").Value = 3790
$ws.Range("C14").Value = 23
$ws.Range("C16").Value = 44
$ws.Range("C17").Value = 83

# --- Literal (non-formula) amount cells ---
$ws.Range("E22").Value = 142494
$ws.Range("E23").Value = 13694

# --- Credit amount ---
$ws.Range("E34").Value = 110000

# --- View / print setup ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.PageSetup.Zoom = 67
